$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DB values (fix to multigraph / updated dbs)
$ws.Range("E8").Value = 1.34
$ws.Range("E9").Value = 1.3

# Move the active selection to E14
$ws.Range("E14").Select()
